# Insert a new data row at row 18 (shifting existing rows 18-53 down to 19-54)
# and populate it with the new weekly price record for "Haba" (Vega Modelo de Temuco).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 18; Excel shifts rows 18..53 down to 19..54
# and the inserted row inherits formatting (e.g. the date format in column D).
$ws.Rows("18").Insert()

# Populate the newly inserted row 18 with the new record's values.
$ws.Range("A18").Value = 10
$ws.Range("B18").Value = 'Vega Modelo de Temuco'
$ws.Range("C18").Value = 'La Araucanía'
$ws.Range("D18").Value = 44519
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = 100112026
$ws.Range("G18").Value = 'Haba'
$ws.Range("H18").Value = 'Sin especificar'
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 55
$ws.Range("K18").Value = 8000
$ws.Range("L18").Value = 9000
$ws.Range("M18").Value = 8455
$ws.Range("N18").Value = '$/saco 25 kilos'
$ws.Range("O18").Value = 'Región del Maule'
$ws.Range("P18").Value = 338
$ws.Range("Q18").Value = 25
$ws.Range("R18").Value = 'Hortaliza'
